# Daily attendance processing - reverse order of "Recorded By" entries
# for cells whose value starts with "System" (column G).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $value = $cell.Value2

    if ($value -ne $null -and $value -is [string] -and $value.StartsWith("System")) {
        $parts = $value.Split(",")
        $trimmed = @()
        foreach ($p in $parts) {
            $trimmed += $p.Trim()
        }

        $reversed = @()
        for ($i = $trimmed.Count - 1; $i -ge 0; $i--) {
            $reversed += $trimmed[$i]
        }

        $newValue = [string]::Join(", ", $reversed)
        $cell.Value = $newValue
    }
}
